# This script applies a rotation of species-observation data among rows 3, 4, 6, 7 and 8
# of the "Artfynd" worksheet. Row 3 and Row 7 swap their content, while rows 4, 6 and 8
# rotate their content (new4 = old6, new6 = old8, new8 = old4). The "Publik kommentar"
# (AC) column follows the same data, including being absent when the target row's
# incoming data had no comment.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 <- old Row 7 data
$ws.Range("A3").Value = 111639167
$ws.Range("B3").Value = 96348
$ws.Range("D3").Value = "VU"
$ws.Range("E3").Value = 220787
$ws.Range("F3").Value = "Knärot"
$ws.Range("G3").Value = "Goodyera repens"
$ws.Range("H3").Value = "(L.) R. Br."
$ws.Range("Q3").Value = 547814.5103353403
$ws.Range("R3").Value = 6926124.461383951
$ws.Range("AC3").Value = "1 planta"

# Row 4 <- old Row 6 data
$ws.Range("A4").Value = 111639174
$ws.Range("Q4").Value = 547803.9854679118
$ws.Range("R4").Value = 6926147.447742103
$ws.Range("AC4").Value = "ca 6 plantor"

# Row 6 <- old Row 8 data
$ws.Range("A6").Value = 111639173
$ws.Range("Q6").Value = 547838.0352795018
$ws.Range("R6").Value = 6926228.915831603
$ws.Range("AC6").Value = "ca 15 plantor"

# Row 7 <- old Row 3 data
$ws.Range("A7").Value = 111639168
$ws.Range("B7").Value = 89686
$ws.Range("D7").Value = "NT"
$ws.Range("E7").Value = 658
$ws.Range("F7").Value = "Rosenticka"
$ws.Range("G7").Value = "Rhodofomes roseus"
$ws.Range("H7").Value = "(Alb. & Schwein.) Kotl. & Pouzar"
$ws.Range("Q7").Value = 548104.1391889038
$ws.Range("R7").Value = 6926477.987023209
$ws.Range("AC7").ClearContents()

# Row 8 <- old Row 4 data
$ws.Range("A8").Value = 111639172
$ws.Range("Q8").Value = 548221.3480213688
$ws.Range("R8").Value = 6926511.607424877
$ws.Range("AC8").ClearContents()
